$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 4) mirroring the existing rows' shape:
# first_name, last_name, email(hyperlink), telephone, password, confirm_password
$ws.Range("A4").Value = "Anisa"
$ws.Range("B4").Value = "Faizi"
$ws.Range("C4").Value = "anisa@gmail.com"
$ws.Range("D4").Value = 2023439873
$ws.Range("E4").Value = "xyz123"
$ws.Range("F4").Value = "xyz123"

# Turn the new email cell into a mailto hyperlink, like C2/C3
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:anisa@gmail.com")

# Hyperlinks.Add re-stamps the cell format; reapply the same Hyperlink
# style the other email cells use so C4 matches C3's styling.
$ws.Range("C4").Style = "Hyperlink"

# Move the active selection, like the workbook shows after the edit
$ws.Range("F7").Select()
